$d = $word.ActiveDocument

# Locate the final paragraph in the document - the one that currently
# just contains the text "asdf" followed by the _GoBack bookmark.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text.TrimEnd([char]13,[char]7) -eq "asdf") {
        $target = $cand
    }
}

$r = $target.Range

# Pull the real WordprocessingML (as a full pkg:package payload) that
# represents this paragraph, so we can splice in replacement markup
# while preserving every namespace declaration, part reference, etc.
# that the runtime expects.
$xml = $r.WordOpenXML

$oldFragment = '<w:r><w:t>asdf</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>'

$newFragment = '<w:proofErr w:type="spellStart"/><w:r><w:t>A</w:t></w:r><w:r><w:t>sdf</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p/><w:p><w:r><w:t>10/2 Things to think about for a successful app</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Offline f</w:t></w:r><w:r><w:t>unctionality</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Personalize the experience</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Offer it free</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Device Orientation</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>'

if ($xml.IndexOf($oldFragment) -lt 0) {
    throw "Could not locate expected 'asdf' run/bookmark markup inside WordOpenXML payload"
}

$newXml = $xml.Replace($oldFragment, $newFragment)

# Strip the synthetic paraId/textId bookkeeping attributes that
# WordOpenXML attaches to the paragraph(s) it scopes over - they are not
# part of the original document and should not leak into the result.
$newXml = $newXml -replace ' w14:paraId="[0-9A-Fa-f]+"', ''
$newXml = $newXml -replace ' w14:textId="[0-9A-Fa-f]+"', ''

$r.InsertXML($newXml)
